# Update time_taken (F column) timestamps on the "data" sheet to reflect the
# later query run captured in this revision.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:22:13.166403"
$data.Range("F3").Value = "2021-10-05 14:22:13.166410"
$data.Range("F4").Value = "2021-10-05 14:22:13.166414"
$data.Range("F5").Value = "2021-10-05 14:22:13.166417"
$data.Range("F6").Value = "2021-10-05 14:22:13.166419"
$data.Range("F7").Value = "2021-10-05 14:22:13.166422"
$data.Range("F8").Value = "2021-10-05 14:22:13.166425"
$data.Range("F9").Value = "2021-10-05 14:22:13.166427"
$data.Range("F10").Value = "2021-10-05 14:22:13.166430"
$data.Range("F11").Value = "2021-10-05 14:22:13.166433"
$data.Range("F12").Value = "2021-10-05 14:22:13.166435"
$data.Range("F13").Value = "2021-10-05 14:22:13.166438"
$data.Range("F14").Value = "2021-10-05 14:22:13.166440"
$data.Range("F15").Value = "2021-10-05 14:22:13.166442"
$data.Range("F16").Value = "2021-10-05 14:22:13.166445"
$data.Range("F17").Value = "2021-10-05 14:22:13.166447"
$data.Range("F18").Value = "2021-10-05 14:22:13.166450"
$data.Range("F19").Value = "2021-10-05 14:22:13.166453"
$data.Range("F20").Value = "2021-10-05 14:22:13.166455"
$data.Range("F21").Value = "2021-10-05 14:22:13.166458"
$data.Range("F22").Value = "2021-10-05 14:22:13.166460"
$data.Range("F23").Value = "2021-10-05 14:22:13.166462"
$data.Range("F24").Value = "2021-10-05 14:22:13.166465"
$data.Range("F25").Value = "2021-10-05 14:22:13.166467"
$data.Range("F26").Value = "2021-10-05 14:22:13.166470"
$data.Range("F27").Value = "2021-10-05 14:22:13.166472"
$data.Range("F28").Value = "2021-10-05 14:22:13.166475"
$data.Range("F29").Value = "2021-10-05 14:22:13.166477"
$data.Range("F30").Value = "2021-10-05 14:22:13.166480"
$data.Range("F31").Value = "2021-10-05 14:22:13.166482"
$data.Range("F32").Value = "2021-10-05 14:22:13.166485"
$data.Range("F33").Value = "2021-10-05 14:22:13.166487"
$data.Range("F34").Value = "2021-10-05 14:22:13.166491"
$data.Range("F35").Value = "2021-10-05 14:22:13.166493"
$data.Range("F36").Value = "2021-10-05 14:22:13.166496"
$data.Range("F37").Value = "2021-10-05 14:22:13.166498"
$data.Range("F38").Value = "2021-10-05 14:22:13.166501"
$data.Range("F39").Value = "2021-10-05 14:22:13.166504"
$data.Range("F40").Value = "2021-10-05 14:22:13.166506"
$data.Range("F41").Value = "2021-10-05 14:22:13.166508"
$data.Range("F42").Value = "2021-10-05 14:22:13.166511"
$data.Range("F43").Value = "2021-10-05 14:22:13.166514"
$data.Range("F44").Value = "2021-10-05 14:22:13.166516"
$data.Range("F45").Value = "2021-10-05 14:22:13.166519"
$data.Range("F46").Value = "2021-10-05 14:22:13.166521"
$data.Range("F47").Value = "2021-10-05 14:22:13.166524"
$data.Range("F48").Value = "2021-10-05 14:22:13.166526"
$data.Range("F49").Value = "2021-10-05 14:22:13.166529"
$data.Range("F50").Value = "2021-10-05 14:22:13.166532"
$data.Range("F51").Value = "2021-10-05 14:22:13.166534"
$data.Range("F52").Value = "2021-10-05 14:22:13.166537"
$data.Range("F53").Value = "2021-10-05 14:22:13.166539"
$data.Range("F54").Value = "2021-10-05 14:22:13.166542"
$data.Range("F55").Value = "2021-10-05 14:22:13.166545"
$data.Range("F56").Value = "2021-10-05 14:22:13.166547"
$data.Range("F57").Value = "2021-10-05 14:22:13.166550"
$data.Range("F58").Value = "2021-10-05 14:22:13.166553"
$data.Range("F59").Value = "2021-10-05 14:22:13.166555"
$data.Range("F60").Value = "2021-10-05 14:22:13.166558"
$data.Range("F61").Value = "2021-10-05 14:22:13.166560"
$data.Range("F62").Value = "2021-10-05 14:22:13.166563"
$data.Range("F63").Value = "2021-10-05 14:22:13.166565"
$data.Range("F64").Value = "2021-10-05 14:22:13.166568"
$data.Range("F65").Value = "2021-10-05 14:22:13.166570"
$data.Range("F66").Value = "2021-10-05 14:22:13.166574"
$data.Range("F67").Value = "2021-10-05 14:22:13.166577"
$data.Range("F68").Value = "2021-10-05 14:22:13.166580"
$data.Range("F69").Value = "2021-10-05 14:22:13.166582"
$data.Range("F70").Value = "2021-10-05 14:22:13.166585"
$data.Range("F71").Value = "2021-10-05 14:22:13.166587"
$data.Range("F72").Value = "2021-10-05 14:22:13.166590"
$data.Range("F73").Value = "2021-10-05 14:22:13.166592"
$data.Range("F74").Value = "2021-10-05 14:22:13.166595"
$data.Range("F75").Value = "2021-10-05 14:22:13.166597"
$data.Range("F76").Value = "2021-10-05 14:22:13.166600"
$data.Range("F77").Value = "2021-10-05 14:22:13.166602"
$data.Range("F78").Value = "2021-10-05 14:22:13.166607"
$data.Range("F79").Value = "2021-10-05 14:22:13.166610"
$data.Range("F80").Value = "2021-10-05 14:22:13.166612"
$data.Range("F81").Value = "2021-10-05 14:22:13.166615"
$data.Range("F82").Value = "2021-10-05 14:22:13.166617"
$data.Range("F83").Value = "2021-10-05 14:22:13.166620"
$data.Range("F84").Value = "2021-10-05 14:22:13.166622"
$data.Range("F85").Value = "2021-10-05 14:22:13.166625"
$data.Range("F86").Value = "2021-10-05 14:22:13.166627"
$data.Range("F87").Value = "2021-10-05 14:22:13.166630"
$data.Range("F88").Value = "2021-10-05 14:22:13.166632"
$data.Range("F89").Value = "2021-10-05 14:22:13.166635"
$data.Range("F90").Value = "2021-10-05 14:22:13.166637"
$data.Range("F91").Value = "2021-10-05 14:22:13.166640"
$data.Range("F92").Value = "2021-10-05 14:22:13.166642"
$data.Range("F93").Value = "2021-10-05 14:22:13.166645"
$data.Range("F94").Value = "2021-10-05 14:22:13.166649"
$data.Range("F95").Value = "2021-10-05 14:22:13.166651"
$data.Range("F96").Value = "2021-10-05 14:22:13.166654"
$data.Range("F97").Value = "2021-10-05 14:22:13.166657"
$data.Range("F98").Value = "2021-10-05 14:22:13.166659"
$data.Range("F99").Value = "2021-10-05 14:22:13.166661"
$data.Range("F100").Value = "2021-10-05 14:22:13.166664"
$data.Range("F101").Value = "2021-10-05 14:22:13.166666"
$data.Range("F102").Value = "2021-10-05 14:22:13.166669"
$data.Range("F103").Value = "2021-10-05 14:22:13.166671"
$data.Range("F104").Value = "2021-10-05 14:22:13.166674"
$data.Range("F105").Value = "2021-10-05 14:22:13.166676"
$data.Range("F106").Value = "2021-10-05 14:22:13.166679"
$data.Range("F107").Value = "2021-10-05 14:22:13.166681"
$data.Range("F108").Value = "2021-10-05 14:22:13.166684"
$data.Range("F109").Value = "2021-10-05 14:22:13.166686"
$data.Range("F110").Value = "2021-10-05 14:22:13.166691"
$data.Range("F111").Value = "2021-10-05 14:22:13.166694"
$data.Range("F112").Value = "2021-10-05 14:22:13.166696"
$data.Range("F113").Value = "2021-10-05 14:22:13.166699"
$data.Range("F114").Value = "2021-10-05 14:22:13.166702"
$data.Range("F115").Value = "2021-10-05 14:22:13.166704"
$data.Range("F116").Value = "2021-10-05 14:22:13.166707"
$data.Range("F117").Value = "2021-10-05 14:22:13.166710"
$data.Range("F118").Value = "2021-10-05 14:22:13.166712"
$data.Range("F119").Value = "2021-10-05 14:22:13.166715"
$data.Range("F120").Value = "2021-10-05 14:22:13.166717"
$data.Range("F121").Value = "2021-10-05 14:22:13.166720"
$data.Range("F122").Value = "2021-10-05 14:22:13.166723"

# Add a new "metadata" sheet (placed after "data") describing the PanelApp
# query that produced this export.
$metadata = $wb.Worksheets.Add($null, $data)
$metadata.Name = "metadata"

$metadata.Range("B1").Value = "data_name"
$metadata.Range("C1").Value = "data_id"
$metadata.Range("D1").Value = "data_version"
$metadata.Range("E1").Value = "data_version_created"
$metadata.Range("F1").Value = "panel_query_time"
$metadata.Range("G1").Value = "panel_get_request"

$metadata.Range("A2").Value = 0
$metadata.Range("B2").Value = "Pigmentary skin disorders"
$metadata.Range("C2").Value = 559
$metadata.Range("D2").Value = "'1.15"
$metadata.Range("E2").Value = "2021-08-31T14:23:49.872379Z"
$metadata.Range("F2").Value = "2021-10-05 14:22:13.163241"
$metadata.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/559/?format=json"

# Match the header styling (bold, bordered, centered) used on the "data"
# sheet's header row / index column by copying its format.
$data.Range("B1").Copy()
$metadata.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)

Write-Output "metadata sheet added"
